# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 46074 (2026-02-21) to 46075 (2026-02-22).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C217").Value = 46075
